# Auto-generated edit script: updates profit-calc sheets with refreshed market-price snapshots.
# Applies per-cell value updates (and a few cell clears / additions) as captured by the
# upstream scheduled-runner diff for Sheets/Behemoth_Profits.xlsx.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 150
$ws.Range("I2").Value = 150
$ws.Range("K2").Value = 150
$ws.Range("M2").Value = -37
$ws.Range("H17").Value = 411.15
$ws.Range("J17").Value = 411.15
$ws.Range("L17").Value = 1233.45
$ws.Range("N17").Value = -1569.45
$ws.Range("H62").Value = 1947.5
$ws.Range("J62").Value = 1913.6666
$ws.Range("L62").Value = 1913.6666
$ws.Range("N62").Value = -3161.6666
$ws.Range("H64").Value = 3879
$ws.Range("I64").Value = 4231.6665
$ws.Range("J64").Value = 3350
$ws.Range("K64").Value = 4231.6665
$ws.Range("L64").Value = 3350
$ws.Range("M64").Value = -3983.6665
$ws.Range("N64").Value = -3846
$ws.Range("H65").Value = 1947.5
$ws.Range("J65").Value = 1913.6666
$ws.Range("L65").Value = 9568.333000000001
$ws.Range("N65").Value = -15808.333
$ws.Range("H67").Value = 3879
$ws.Range("I67").Value = 4231.6665
$ws.Range("J67").Value = 3350
$ws.Range("K67").Value = 4231.6665
$ws.Range("L67").Value = 3350
$ws.Range("M67").Value = -3373.6665
$ws.Range("N67").Value = -5066
$ws.Range("H112").Value = 2699.75
$ws.Range("J112").Value = 2666.389
$ws.Range("L112").Value = 7999.167
$ws.Range("N112").Value = -10215.167
$ws.Range("H132").Value = 3231.3684
$ws.Range("I132").Value = 2846.8235
$ws.Range("K132").Value = 8540.470499999999
$ws.Range("M132").Value = -6010.470499999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 69335.39999999999
$ws.Range("J31").Value = 111226
$ws.Range("L31").Value = 111226
$ws.Range("N31").Value = -111814
$ws.Range("H32").Value = 26327040
$ws.Range("I32").Value = 26327040
$ws.Range("K32").Value = 26327040
$ws.Range("M32").Value = -26326753
$ws.Range("H74").Value = 14716325
$ws.Range("I74").Value = 20834432
$ws.Range("K74").Value = 20834432
$ws.Range("M74").Value = -20833558
$ws.Range("H77").Value = 14716325
$ws.Range("I77").Value = 20834432
$ws.Range("K77").Value = 104172160
$ws.Range("M77").Value = -104167792
$ws.Range("H124").Value = 47474.668
$ws.Range("J124").Value = 47474.668
$ws.Range("L124").Value = 47474.668
$ws.Range("N124").Value = -57294.668
$ws.Range("H125").Value = 70712.5
$ws.Range("J125").Value = 70712.5
$ws.Range("L125").Value = 70712.5
$ws.Range("N125").Value = -80552.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H96").Value = 36173.332
$ws.Range("J96").Value = 70896.60000000001
$ws.Range("L96").Value = 70896.60000000001
$ws.Range("N96").Value = -76388.60000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 865993.1
$ws.Range("J31").Value = 2123013.8
$ws.Range("L31").Value = 2123013.8
$ws.Range("N31").Value = -2123603.8
$ws.Range("H34").Value = 865993.1
$ws.Range("J34").Value = 2123013.8
$ws.Range("L34").Value = 2123013.8
$ws.Range("N34").Value = -2123417.8
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("N99").ClearContents()
$ws.Range("H107").Value = 1547.6364
$ws.Range("I107").Value = 1137.25
$ws.Range("K107").Value = 1137.25
$ws.Range("M107").Value = 782.75
$ws.Range("H122").Value = 4699
$ws.Range("I122").Value = 4699
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 14097
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -11647
$ws.Range("N122").ClearContents()
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").ClearContents()
$ws.Range("H134").Value = 202661.64
$ws.Range("I134").Value = 264694.28
$ws.Range("K134").Value = 794082.8400000001
$ws.Range("M134").Value = -791547.8400000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 44531.223
$ws.Range("J2").Value = 100045.75
$ws.Range("L2").Value = 600274.5
$ws.Range("N2").Value = -600500.5
$ws.Range("H11").Value = 344.0345
$ws.Range("I11").Value = 274.30435
$ws.Range("J11").Value = 611.3333
$ws.Range("K11").Value = 822.91305
$ws.Range("L11").Value = 1833.9999
$ws.Range("M11").Value = -682.91305
$ws.Range("N11").Value = -2113.9999
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("N23").ClearContents()
$ws.Range("H37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").ClearContents()
$ws.Range("H86").Value = 144.375
$ws.Range("H87").Value = 20199.8
$ws.Range("I87").Value = 18999.75
$ws.Range("K87").Value = 56999.25
$ws.Range("M87").Value = -55751.25
$ws.Range("H89").Value = 144.375
$ws.Range("H90").Value = 20199.8
$ws.Range("I90").Value = 18999.75
$ws.Range("K90").Value = 170997.75
$ws.Range("M90").Value = -164757.75
$ws.Range("H92").Value = 2001884.2
$ws.Range("I92").Value = 2501730.8
$ws.Range("K92").Value = 7505192.399999999
$ws.Range("M92").Value = -7503944.399999999
$ws.Range("H128").Value = 419995
$ws.Range("I128").Value = 419995
$ws.Range("K128").Value = 1259985
$ws.Range("M128").Value = -1255005
$ws.Range("H140").Value = 502263.5
$ws.Range("I140").Value = 502263.5
$ws.Range("K140").Value = 1506790.5
$ws.Range("M140").Value = -1501610.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 3400.25
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()
$ws.Range("H62").Value = 100437
$ws.Range("J62").Value = 100437
$ws.Range("L62").Value = 100437
$ws.Range("N62").Value = -101809
$ws.Range("H65").Value = 100437
$ws.Range("J65").Value = 100437
$ws.Range("L65").Value = 301311
$ws.Range("N65").Value = -308175
$ws.Range("H70").Value = 3200
$ws.Range("I70").Value = 2700
$ws.Range("J70").Value = 3700
$ws.Range("K70").Value = 2700
$ws.Range("L70").Value = 3700
$ws.Range("M70").Value = -2430
$ws.Range("N70").Value = -4240
$ws.Range("H73").Value = 3200
$ws.Range("I73").Value = 2700
$ws.Range("J73").Value = 3700
$ws.Range("K73").Value = 2700
$ws.Range("L73").Value = 3700
$ws.Range("M73").Value = -1764
$ws.Range("N73").Value = -5572
$ws.Range("H97").Value = 787.96155
$ws.Range("I97").Value = 718.381
$ws.Range("J97").Value = 1080.2
$ws.Range("K97").Value = 718.381
$ws.Range("L97").Value = 1080.2
$ws.Range("M97").Value = -222.381
$ws.Range("N97").Value = -2072.2
$ws.Range("H109").Value = 47500
$ws.Range("J109").Value = 47500
$ws.Range("L109").Value = 47500
$ws.Range("N109").Value = -49580
$ws.Range("H132").Value = 58826428
$ws.Range("I132").Value = 66669550
$ws.Range("K132").Value = 200008650
$ws.Range("M132").Value = -200006120

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 1200
$ws.Range("I2").Value = 1200
$ws.Range("K2").Value = 1200
$ws.Range("M2").Value = -1088
$ws.Range("H127").Value = 140000
$ws.Range("J127").Value = 140000
$ws.Range("L127").Value = 140000
$ws.Range("N127").Value = -149920
$ws.Range("H136").Value = 96045.83
$ws.Range("I136").Value = 58391.668
$ws.Range("K136").Value = 175175.004
$ws.Range("M136").Value = -172625.004

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 10766536
$ws.Range("J75").Value = 10766536
$ws.Range("L75").Value = 10766536
$ws.Range("N75").Value = -10768408
$ws.Range("H78").Value = 10766536
$ws.Range("J78").Value = 10766536
$ws.Range("L78").Value = 32299608
$ws.Range("N78").Value = -32308968
$ws.Range("H107").Value = 27778856
$ws.Range("J107").Value = 690
$ws.Range("L107").Value = 2070
$ws.Range("N107").Value = -5910
$ws.Range("H122").Value = 3579.6843
$ws.Range("I122").Value = 2357.6
$ws.Range("J122").Value = 5929.846
$ws.Range("K122").Value = 7072.799999999999
$ws.Range("L122").Value = 17789.538
$ws.Range("M122").Value = -4622.799999999999
$ws.Range("N122").Value = -22689.538

Write-Output "Applied 219 cell updates across 8 sheets."
